$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1=14, Q1=15 using the same formatting as O1 ---
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value = 15
$ws.Application.CutCopyMode = $false

# --- Data rows 2-25 ---
# Columns I, K, M, O swap their 1/2 value; new columns P and Q are appended
# with value 2 (no special style, matching the plain data cells).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column
}
